# Add solutions/documentation rows for LeetCode problems:
#   3583 - Count Special Triplets
#   3577 - Count the Number of Computer Unlocking Permutations
#   1874 - Minimize Product Sum of Two Arrays
# plus a trailing "date stamp" row, mirroring the existing sheet pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Clone the formatting (styles) of the last existing data row (245) down
# --- into the three new rows, so columns keep their usual per-column style
# --- (A/D/E/F/G centered, B/C wrap-text, H/I date-formatted) without
# --- introducing any new style/numFmt definitions.
$ws.Range("A245:I245").Copy()
$ws.Range("A246:I248").PasteSpecial($xlPasteFormats)

# --- Row 246: 3583. Count Special Triplets
$ws.Range("A246").Value = 3583
$ws.Range("B246").Value = "Count Special Triplets"
$ws.Range("C246").Value = "#array #hash-table "
$ws.Range("D246").Value = "medium"
$ws.Range("E246").Value = 0
$ws.Range("F246").Value = 1
$ws.Range("G246").Value = 30
$ws.Range("H246").Value = 46000
$ws.Range("I246").Value = 46000

# --- Row 247: 3577. Count the Number of Computer Unlocking Permutations
$ws.Range("A247").Value = 3577
$ws.Range("B247").Value = "Count the Number of Computer Unlocking Permutations"
$ws.Range("C247").Value = "#math #array"
$ws.Range("D247").Value = "medium"
$ws.Range("E247").Value = 1
$ws.Range("F247").Value = 0
$ws.Range("G247").Value = 10
$ws.Range("H247").Value = 46001
$ws.Range("I247").Value = 46001

# --- Row 248: 1874. Minimize Product Sum of Two Arrays
$ws.Range("A248").Value = 1874
$ws.Range("B248").Value = "Minimize Product Sum of Two Arrays"
$ws.Range("C248").Value = "#math #array "
$ws.Range("D248").Value = "medium"
$ws.Range("E248").Value = 1
$ws.Range("F248").Value = 0
$ws.Range("G248").Value = 4
$ws.Range("H248").Value = 46001
$ws.Range("I248").Value = 46001

# --- Row 249: trailing date-stamp only row (H/I), matching the workbook's
# --- existing convention of blank "Last Update" marker rows.
$ws.Range("H245:I245").Copy()
$ws.Range("H249:I249").PasteSpecial($xlPasteFormats)
$ws.Range("H249").Value = 46001
$ws.Range("I249").Value = 46001

$excel.CutCopyMode = $false

# --- Row heights follow Excel's wrap-text auto-height (17pt per wrapped
# --- line) for the rows whose Name/Tags text wraps across columns B/C.
$ws.Rows.Item(246).RowHeight = 17
$ws.Rows.Item(247).RowHeight = 51
$ws.Rows.Item(248).RowHeight = 34

# --- Leave the same selection state the author's session ended in.
$null = $ws.Range("H249:I249").Select()
